# Atualizado por script em 05-11-2023 08:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data (columns F..V) between row 5 and row 6 ---
# (Indice/pais/torneio/temporada/data_partida in columns A-E stay put;
#  only the match details moved rows.)
for ($col = 6; $col -le 22; $col++) {
    $v5 = $ws.Cells.Item(5, $col).Value()
    $v6 = $ws.Cells.Item(6, $col).Value()
    $ws.Cells.Item(5, $col).Value = $v6
    $ws.Cells.Item(6, $col).Value = $v5
}

# --- Append two new match rows (43 and 44) ---

# Row 43: Birkirkara 0-2 Marsaxlokk
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(43, 1).PasteSpecial(-4122)
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "malta"
$ws.Cells.Item(43, 3).Value = "premier-league"
$ws.Cells.Item(43, 4).Value = "2023-2024"
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(43, 5).PasteSpecial(-4122)
$ws.Cells.Item(43, 5).Value = 45234.58333333334
$ws.Cells.Item(43, 6).Value = "Birkirkara"
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = "Marsaxlokk"
$ws.Cells.Item(43, 9).Value = 2
$ws.Cells.Item(43, 10).Value = 1.75
$ws.Cells.Item(43, 11).Value = "03/11/2023 08:42"
$ws.Cells.Item(43, 12).Value = 1.89
$ws.Cells.Item(43, 13).Value = "04/11/2023 13:58"
$ws.Cells.Item(43, 14).Value = 3.41
$ws.Cells.Item(43, 15).Value = "03/11/2023 08:42"
$ws.Cells.Item(43, 16).Value = 2.95
$ws.Cells.Item(43, 17).Value = "04/11/2023 13:58"
$ws.Cells.Item(43, 18).Value = 4.12
$ws.Cells.Item(43, 19).Value = "03/11/2023 08:42"
$ws.Cells.Item(43, 20).Value = 4.89
$ws.Cells.Item(43, 21).Value = "04/11/2023 13:58"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/birkirkara-marsaxlokk/vDEV123d/"

# Row 44: Gzira 1-1 Valletta
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(44, 1).PasteSpecial(-4122)
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "malta"
$ws.Cells.Item(44, 3).Value = "premier-league"
$ws.Cells.Item(44, 4).Value = "2023-2024"
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(44, 5).PasteSpecial(-4122)
$ws.Cells.Item(44, 5).Value = 45234.6875
$ws.Cells.Item(44, 6).Value = "Gzira"
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = "Valletta"
$ws.Cells.Item(44, 9).Value = 1
$ws.Cells.Item(44, 10).Value = 1.93
$ws.Cells.Item(44, 11).Value = "03/11/2023 08:42"
$ws.Cells.Item(44, 12).Value = 2.45
$ws.Cells.Item(44, 13).Value = "04/11/2023 16:20"
$ws.Cells.Item(44, 14).Value = 3.21
$ws.Cells.Item(44, 15).Value = "03/11/2023 08:42"
$ws.Cells.Item(44, 16).Value = 2.85
$ws.Cells.Item(44, 17).Value = "04/11/2023 16:20"
$ws.Cells.Item(44, 18).Value = 3.44
$ws.Cells.Item(44, 19).Value = "03/11/2023 08:42"
$ws.Cells.Item(44, 20).Value = 3.21
$ws.Cells.Item(44, 21).Value = "04/11/2023 16:20"
$ws.Cells.Item(44, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/gzira-valletta/tIhEQotc/"

"done"
